$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1889
$ws.Range("I38").Value = 71
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 213
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = 159
$ws.Range("N38").Value = -9744

$ws.Range("H39").Value = 96.36364
$ws.Range("I39").Value = 51.11111
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 153.33333
$ws.Range("L39").Value = 900
$ws.Range("M39").Value = 142.66667
$ws.Range("N39").Value = -1492

$ws.Range("H58").Value = 835
$ws.Range("I58").Value = 521.4706
$ws.Range("K58").Value = 1564.4118
$ws.Range("M58").Value = -1414.4118

$ws.Range("H87").Value = 25962.375
$ws.Range("J87").Value = 25962.375
$ws.Range("L87").Value = 25962.375
$ws.Range("N87").Value = -28458.375

$ws.Range("H90").Value = 25962.375
$ws.Range("J90").Value = 25962.375
$ws.Range("L90").Value = 77887.125
$ws.Range("N90").Value = -90367.125

$ws.Range("H111").Value = 1461.25
$ws.Range("I111").Value = 1611.7693
$ws.Range("J111").Value = 1181.7142
$ws.Range("K111").Value = 4835.3079
$ws.Range("L111").Value = 3545.1426
$ws.Range("M111").Value = -1768.3079
$ws.Range("N111").Value = -9679.142599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 3751.5
$ws.Range("I16").Value = 2335.3333
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 2335.3333
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -2048.3333
$ws.Range("N16").Value = -8574

$ws.Range("H88").Value = 125101660
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 125101660
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 125101660
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -125102472

$ws.Range("H91").Value = 125101660
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 125101660
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 125101660
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -125104468

$ws.Range("H122").Value = 2745.946
$ws.Range("I122").Value = 2364.6667
$ws.Range("J122").Value = 3775.4
$ws.Range("K122").Value = 7094.000100000001
$ws.Range("L122").Value = 11326.2
$ws.Range("M122").Value = -4644.000100000001
$ws.Range("N122").Value = -16226.2

$ws.Range("H132").Value = 34785.87
$ws.Range("I132").Value = 49583.76
$ws.Range("J132").Value = 3710.3
$ws.Range("K132").Value = 148751.28
$ws.Range("L132").Value = 11130.9
$ws.Range("M132").Value = -146221.28
$ws.Range("N132").Value = -16190.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18002
$ws.Range("I82").Value = 4324.7
$ws.Range("J82").Value = 28523
$ws.Range("K82").Value = 4324.7
$ws.Range("L82").Value = 28523
$ws.Range("M82").Value = -3941.7
$ws.Range("N82").Value = -29289

$ws.Range("H85").Value = 18002
$ws.Range("I85").Value = 4324.7
$ws.Range("J85").Value = 28523
$ws.Range("K85").Value = 4324.7
$ws.Range("L85").Value = 28523
$ws.Range("M85").Value = -2998.7
$ws.Range("N85").Value = -31175

$ws.Range("H94").Value = 10184.242
$ws.Range("I94").Value = 5428.5454
$ws.Range("K94").Value = 5428.5454
$ws.Range("M94").Value = -4977.5454

$ws.Range("H99").Value = 1760.0646
$ws.Range("I99").Value = 1511.5264
$ws.Range("J99").Value = 2153.5833
$ws.Range("K99").Value = 1511.5264
$ws.Range("L99").Value = 2153.5833
$ws.Range("M99").Value = -13.52639999999997
$ws.Range("N99").Value = -5149.5833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 7000
$ws.Range("J15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("N15").Value = -7340

$ws.Range("H31").Value = 2049.2554
$ws.Range("I31").Value = 1261.4814
$ws.Range("J31").Value = 3112.75
$ws.Range("K31").Value = 1261.4814
$ws.Range("L31").Value = 3112.75
$ws.Range("M31").Value = -966.4813999999999
$ws.Range("N31").Value = -3702.75

$ws.Range("H34").Value = 2049.2554
$ws.Range("I34").Value = 1261.4814
$ws.Range("J34").Value = 3112.75
$ws.Range("K34").Value = 1261.4814
$ws.Range("L34").Value = 3112.75
$ws.Range("M34").Value = -1059.4814
$ws.Range("N34").Value = -3516.75

$ws.Range("H132").Value = 1808.742
$ws.Range("I132").Value = 1581.0454
$ws.Range("K132").Value = 4743.1362
$ws.Range("M132").Value = -2213.1362

$ws.Range("H134").Value = 3826.9768
$ws.Range("I134").Value = 3938.7778
$ws.Range("J134").Value = 3252
$ws.Range("K134").Value = 11816.3334
$ws.Range("L134").Value = 9756
$ws.Range("M134").Value = -9281.3334
$ws.Range("N134").Value = -14826

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 939.2759
$ws.Range("J5").Value = 1335.3572
$ws.Range("L5").Value = 4006.0716
$ws.Range("N5").Value = -4230.071599999999

$ws.Range("H113").Value = 1390.525
$ws.Range("I113").Value = 526.6316
$ws.Range("J113").Value = 2172.1428
$ws.Range("K113").Value = 1579.8948
$ws.Range("L113").Value = 6516.428400000001
$ws.Range("M113").Value = 590.1052
$ws.Range("N113").Value = -10856.4284

$ws.Range("H121").Value = 1297.4166
$ws.Range("I121").Value = 349.8
$ws.Range("J121").Value = 1974.2858
$ws.Range("K121").Value = 1049.4
$ws.Range("L121").Value = 5922.857400000001
$ws.Range("M121").Value = 260.5999999999999
$ws.Range("N121").Value = -8542.857400000001

$ws.Range("H135").Value = 939.2759
$ws.Range("J135").Value = 1335.3572
$ws.Range("L135").Value = 12018.2148
$ws.Range("N135").Value = -17088.2148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 12407.333
$ws.Range("I31").Value = 7500
$ws.Range("J31").Value = 22222
$ws.Range("K31").Value = 7500
$ws.Range("L31").Value = 22222
$ws.Range("M31").Value = -7208
$ws.Range("N31").Value = -22806

$ws.Range("H37").Value = 12407.333
$ws.Range("I37").Value = 7500
$ws.Range("J37").Value = 22222
$ws.Range("K37").Value = 7500
$ws.Range("L37").Value = 22222
$ws.Range("M37").Value = -7223
$ws.Range("N37").Value = -22776

$ws.Range("H70").Value = 3683970.5
$ws.Range("I70").Value = 7356691
$ws.Range("J70").Value = 11250
$ws.Range("K70").Value = 7356691
$ws.Range("L70").Value = 11250
$ws.Range("M70").Value = -7356421
$ws.Range("N70").Value = -11790

$ws.Range("H73").Value = 3683970.5
$ws.Range("I73").Value = 7356691
$ws.Range("J73").Value = 11250
$ws.Range("K73").Value = 7356691
$ws.Range("L73").Value = 11250
$ws.Range("M73").Value = -7355755
$ws.Range("N73").Value = -13122

$ws.Range("H122").Value = 1860.8077
$ws.Range("I122").Value = 1778.3158
$ws.Range("J122").Value = 2084.7144
$ws.Range("K122").Value = 5334.9474
$ws.Range("L122").Value = 6254.1432
$ws.Range("M122").Value = -2884.9474
$ws.Range("N122").Value = -11154.1432

$ws.Range("H126").Value = 2264.7058
$ws.Range("I126").Value = 1875
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 5625
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -3155
$ws.Range("N126").Value = -14540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1664.0769
$ws.Range("I40").Value = 1382.8
$ws.Range("K40").Value = 1382.8
$ws.Range("M40").Value = -1246.8

$ws.Range("H133").Value = 49326
$ws.Range("J133").Value = 49326
$ws.Range("L133").Value = 49326
$ws.Range("N133").Value = -54386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 878.8333
$ws.Range("I96").Value = 851.5
$ws.Range("J96").Value = 892.5
$ws.Range("K96").Value = 851.5
$ws.Range("L96").Value = 892.5
$ws.Range("M96").Value = 521.5
$ws.Range("N96").Value = -3638.5

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H132").Value = 3583.5957
$ws.Range("I132").Value = 4127.533
$ws.Range("J132").Value = 2623.7058
$ws.Range("K132").Value = 12382.599
$ws.Range("L132").Value = 7871.117400000001
$ws.Range("M132").Value = -9852.599000000002
$ws.Range("N132").Value = -12931.1174
